$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: "Density (lb/in^3)" label in column A, with value 0.284 across B:Y
$ws.Range("A11").Value = "Density (lb/in^3)"
$ws.Range("A11").Style = "Good"

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")
foreach ($col in $cols) {
    $ws.Range("$col`11").Value = 0.28399999999999997
}
